# Adds a new leading "index" column (A) to Sheet1, shifting the existing
# image_uuid / image_name / image_tags columns right to B:D, and fills the
# new column with a centered, bordered 1-based row counter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122
# xlCenter
$xlCenter = -4108

# 1) Insert a new blank column before column A; this shifts the existing
#    data (incl. formatting, column widths, dimension) from A:C to B:D.
$ws.Columns.Item(1).Insert()

# 2) Header cell A1: copy the header formatting from B1, then set the text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial($xlPasteFormats)
$ws.Range("A1").Value() = "index"

# 3) Data cells A2:A18: copy the bordered formatting from column B, then
#    center the values and fill in sequential index numbers (1..17).
$ws.Range("B2:B18").Copy()
$ws.Range("A2:A18").PasteSpecial($xlPasteFormats)
$ws.Range("A2:A18").HorizontalAlignment = $xlCenter

for ($i = 2; $i -le 18; $i++) {
    $ws.Cells.Item($i, 1).Value() = ($i - 1)
}

# 4) Update the active selection to match the edited workbook.
$ws.Range("B21").Select() | Out-Null

Write-Output "applied index column edit"
